{"js": "const paras = context.document.body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nconst titlePara = paras.items[0];\ntitlePara.insertText(\"Week 6 Reading Guide Part 2: Sampling Variability\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$d.Content.Text.Substring(0, 80)\n"}
